# Update tracker data: remove the "G1 / Test1" row (original row 2).
# Deleting the entire row shifts all subsequent rows up by one,
# turning former rows 3-7 (G2..G6) into the new rows 2-6, and
# shrinking the used range from A1:F7 to A1:F6 — matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

$wb.Save()
